$d = $word.ActiveDocument

# ==================================================================
# Hunk 1 + 2 combined: fix "ValidateFormDate.jsp" -> "ValidateFormData.jsp"
# typo, bold the filename, and relocate the "_GoBack" bookmark from
# its old spot (between "ValidateFormData" and ".jsp" further down
# in the document) to the newly-bolded filename here (between the
# "a" and ".jsp").
# ==================================================================

$target = $d.Content
$target.Find.Execute("ValidateFormDate.jsp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fnStart = $target.Start
$fnEnd = $target.End

# Guard: the run right after ("the following fi...") happens to share
# identical run formatting with the text we are about to edit. Any
# text-content mutation in this engine re-normalises (merges)
# adjacent same-formatted runs across the whole paragraph, which
# would wrongly fuse that run into ours. Bolding it first gives it a
# distinct format so it survives untouched; we revert the bold once
# our text edits are done.
$guard = $d.Content
$guard.Find.Execute("the following fi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$guard.Bold = $true

# Fix the typo: "ValidateFormD­ate" -> "ValidateFormData" (just the
# single "e" -> "a").
$typo = $d.Range($fnStart + 15, $fnStart + 16)
$typo.Text = "a"

# Bold the whole (now-fixed) filename "ValidateFormData.jsp".
$filename = $d.Range($fnStart, $fnEnd)
$filename.Bold = $true

# Release the guard.
$guard2 = $d.Content
$guard2.Find.Execute("the following fi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$guard2.Bold = $false

# Split "ValidateFormDat" from "a" into separate runs: add then
# immediately delete a throwaway bookmark at that point purely to
# force the run boundary (a text/bookmark structural edit splits a
# run even when both sides share identical formatting; deleting the
# bookmark afterwards leaves that split in place).
$splitPoint = $d.Range($fnStart + 15, $fnStart + 15)
$d.Bookmarks.Add("zzsplit", $splitPoint) | Out-Null
$d.Bookmarks("zzsplit").Delete()

# Move "_GoBack": drop it from its old home further down the document
# (between "ValidateFormData" and ".jsp") ...
$d.Bookmarks("_GoBack").Delete()

# ... and re-create it here, between the fixed "a" and ".jsp".
$bmPoint = $d.Range($fnStart + 16, $fnStart + 16)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ==================================================================
# Hunk 3: bold "ProcessCustomerDataRequest.jsp" in
# `forward the user to "ProcessCustomerDataRequest.jsp".` (first
# occurrence only -- the later, already-bold occurrence is untouched).
# ==================================================================

$pcdr = $d.Content
$pcdr.Find.Execute("ProcessCustomerDataRequest.jsp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pcdr.Bold = $true
